$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr2 = $sec.Headers.Item(2)
$shp = $hdr2.Range.InlineShapes.Item(1)
$shpRange = $shp.Range
$shp2 = $shpRange.InlineShapes.Item(1)
$shp2.AlternativeText = "TESTDESC"
Write-Output "done"
